$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cars")

# Remove the two "~UC_Sets: Ts_S: " header rows (row 1 of each UC block).
# Delete from the bottom up so earlier row numbers remain valid.
$ws.Rows.Item(8).Delete() | Out-Null
$ws.Rows.Item(1).Delete() | Out-Null
